# Apply "Todays Changes To Flows and Control Relationships"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Claim Filing")

# Add new values to column H, matching the flow/control relationship additions
$ws.Range("H4").Value = "What if NMR is NOT PC ?"
$ws.Range("H8").Value = "What if BPTW is NOT PC ?"

# Set width for the new column H
$ws.Columns.Item(8).ColumnWidth = 27.3984375

# Update the active selection to reflect the new last-used cell
$ws.Range("G14").Select()
